$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The weekly update shifts the date/volume/price data for rows 2-5
# down by one row (with wraparound): row2 <- old row5, row3 <- old row2,
# row4 <- old row3, row5 <- old row4. Columns affected: D, M, N, O, P, S.

$cols = @("D", "M", "N", "O", "P", "S")

# Capture the original (before-edit) values for rows 2-5 first,
# since we will be overwriting them in place. Use Value2 for reads
# (plain numeric, avoids any date-variant formatting surprises).
$orig = @{}
foreach ($r in 2..5) {
    foreach ($col in $cols) {
        $orig["$col$r"] = $ws.Range("$col$r").Value2
    }
}

# New row index -> source (old) row index
$rowMap = @{ 2 = 5; 3 = 2; 4 = 3; 5 = 4 }

foreach ($newRow in 2..5) {
    $oldRow = $rowMap[$newRow]
    foreach ($col in $cols) {
        $ws.Range("$col$newRow").Value2 = $orig["$col$oldRow"]
    }
}
